$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1) "ISIC to BLS Map" sheet: split "ISIC 05T06" (row 3) into two
#    rows - "ISIC 05" (Coal mining) and "ISIC 06" (Oil and gas
#    extraction). Insert a new row above the existing row 3, copy
#    the mapped-category/lookup formula down from the old row 3,
#    then relabel the old row (now row 4).
# ---------------------------------------------------------------
$map = $wb.Worksheets.Item("ISIC to BLS Map")

# Insert a new blank row before row 3 (old row 3 shifts to row 4,
# formulas and relative references shift automatically).
$map.Rows.Item(3).Insert()

# Copy formatting/formulas from the row that used to be row 3 (now
# row 4) into the freshly inserted row 3, so C3/D3 keep the same
# "Mining, quarrying, and oil and gas extraction" lookup as before.
$map.Rows.Item(4).Copy()
$map.Rows.Item(3).PasteSpecial(-4104) | Out-Null
$excel.CutCopyMode = 0

# New row 3: ISIC 05 / Coal mining
$map.Range("A3").Value = "ISIC 05"
$map.Range("B3").Value = "Coal mining"

# Old row (now row 4): relabel to ISIC 06 / Oil and gas extraction
$map.Range("A4").Value = "ISIC 06"
$map.Range("B4").Value = "Oil and gas extraction"

# ---------------------------------------------------------------
# 2) "URPbIC" sheet: mirror the split with a new column C, so the
#    ISIC codes line up with the updated "ISIC to BLS Map" table.
# ---------------------------------------------------------------
$urp = $wb.Worksheets.Item("URPbIC")

# Insert a new blank column before column C (old column C shifts to
# column D).
$urp.Columns.Item(3).Insert()

# Copy formatting/formulas from the column that used to be column C
# (now column D) into the freshly inserted column C.
$urp.Columns.Item(4).Copy()
$urp.Columns.Item(3).PasteSpecial(-4104) | Out-Null
$excel.CutCopyMode = 0

# New column header: ISIC 05
$urp.Range("C1").Value = "ISIC 05"
# Shifted column header (was "ISIC 05T06"): ISIC 06
$urp.Range("D1").Value = "ISIC 06"

# ---------------------------------------------------------------
# 3) Workbook calculation settings: switch on iterative calculation
#    (matches the updated calcPr in workbook.xml).
# ---------------------------------------------------------------
$excel.Iteration = $true
$excel.MaxChange = 0.000010000000000000001

$wb.RefreshAll()
$excel.CalculateFullRebuild()
